$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44350
$ws.Range("L2").Value = 'Especial'
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 24000
$ws.Range("O2").Value = 24000
$ws.Range("P2").Value = 24000
$ws.Range("R2").Value = 'Provincia de Limarí'
$ws.Range("S2").Value = 1333

# Row 3
$ws.Range("D3").Value = 44298
$ws.Range("M3").Value = 160
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 20000
$ws.Range("S3").Value = 1111

# Row 4
$ws.Range("D4").Value = 44271
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("S4").Value = 833

# Row 5
$ws.Range("D5").Value = 44258
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 100
$ws.Range("R5").Value = 'Provincia de Limarí'

# Row 6
$ws.Range("D6").Value = 44284
$ws.Range("L6").Value = 'Especial'
$ws.Range("M6").Value = 120
$ws.Range("N6").Value = 13000
$ws.Range("O6").Value = 13000
$ws.Range("P6").Value = 13000
$ws.Range("S6").Value = 722

# Row 7
$ws.Range("D7").Value = 44284
$ws.Range("L7").Value = 'Extra (doble especial)'
$ws.Range("M7").Value = 100

# Row 8
$ws.Range("D8").Value = 44284
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 12000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 12000
$ws.Range("S8").Value = 667

# Row 9
$ws.Range("D9").Value = 44330
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 23000
$ws.Range("O9").Value = 23000
$ws.Range("P9").Value = 23000
$ws.Range("S9").Value = 1278

# Row 10
$ws.Range("D10").Value = 44299
$ws.Range("L10").Value = 'Especial'
$ws.Range("M10").Value = 170
$ws.Range("N10").Value = 18000
$ws.Range("O10").Value = 18000
$ws.Range("P10").Value = 18000
$ws.Range("R10").Value = 'Provincia de Melipilla'
$ws.Range("S10").Value = 1000

# Row 11
$ws.Range("D11").Value = 44299
$ws.Range("N11").Value = 16000
$ws.Range("O11").Value = 16000
$ws.Range("P11").Value = 16000
$ws.Range("Q11").Value = '$/caja 18 kilos'
$ws.Range("R11").Value = 'Provincia de Melipilla'
$ws.Range("S11").Value = 889
$ws.Range("T11").Value = 18

# Row 12
$ws.Range("D12").Value = 44252
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 140
$ws.Range("N12").Value = 13000
$ws.Range("O12").Value = 13000
$ws.Range("P12").Value = 13000
$ws.Range("Q12").Value = '$/caja 18 kilos'
$ws.Range("S12").Value = 722
$ws.Range("T12").Value = 18

# Row 13
$ws.Range("D13").Value = 44267
$ws.Range("M13").Value = 120
$ws.Range("Q13").Value = '$/caja 18 kilos'
$ws.Range("S13").Value = 722
$ws.Range("T13").Value = 18

# Row 14
$ws.Range("D14").Value = 44292
$ws.Range("L14").Value = 'Especial'
$ws.Range("M14").Value = 150
$ws.Range("N14").Value = 16000
$ws.Range("O14").Value = 16000
$ws.Range("P14").Value = 16000
$ws.Range("Q14").Value = '$/caja 18 kilos'
$ws.Range("S14").Value = 889
$ws.Range("T14").Value = 18

# Row 15
$ws.Range("D15").Value = 44292
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 80
$ws.Range("N15").Value = 14000
$ws.Range("O15").Value = 14000
$ws.Range("P15").Value = 14000
$ws.Range("S15").Value = 778

# Row 16
$ws.Range("D16").Value = 44300
$ws.Range("L16").Value = 'Especial'
$ws.Range("M16").Value = 120
$ws.Range("N16").Value = 18000
$ws.Range("O16").Value = 18000
$ws.Range("P16").Value = 18000
$ws.Range("S16").Value = 1000

# Row 17
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 16000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 16000
$ws.Range("S17").Value = 889

# Row 18
$ws.Range("D18").Value = 44277
$ws.Range("L18").Value = 'Especial'
$ws.Range("M18").Value = 200
$ws.Range("N18").Value = 15000
$ws.Range("O18").Value = 15000
$ws.Range("P18").Value = 15000
$ws.Range("R18").Value = 'Provincia de Limarí'
$ws.Range("S18").Value = 833

# Row 19
$ws.Range("D19").Value = 44224
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 120
$ws.Range("N19").Value = 18000
$ws.Range("O19").Value = 18000
$ws.Range("P19").Value = 18000
$ws.Range("R19").Value = 'Provincia de Limarí'
$ws.Range("S19").Value = 1125

# Row 20
$ws.Range("D20").Value = 44295
$ws.Range("L20").Value = 'Segunda'
$ws.Range("N20").Value = 10000
$ws.Range("O20").Value = 10000
$ws.Range("P20").Value = 10000
$ws.Range("Q20").Value = '$/caja 18 kilos'
$ws.Range("S20").Value = 556
$ws.Range("T20").Value = 18

# Row 21
$ws.Range("D21").Value = 44309
$ws.Range("L21").Value = 'Especial'
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 20000
$ws.Range("O21").Value = 20000
$ws.Range("P21").Value = 20000
$ws.Range("Q21").Value = '$/caja 18 kilos'
$ws.Range("R21").Value = 'Provincia de Melipilla'
$ws.Range("S21").Value = 1111
$ws.Range("T21").Value = 18

# Row 22
$ws.Range("D22").Value = 44309
$ws.Range("M22").Value = 60
$ws.Range("N22").Value = 18000
$ws.Range("O22").Value = 18000
$ws.Range("P22").Value = 18000
$ws.Range("S22").Value = 1000

# Row 23
$ws.Range("D23").Value = 44274
$ws.Range("M23").Value = 200
$ws.Range("N23").Value = 14000
$ws.Range("O23").Value = 14000
$ws.Range("P23").Value = 14000
$ws.Range("Q23").Value = '$/caja 16 kilos'
$ws.Range("R23").Value = 'Provincia de Melipilla'
$ws.Range("S23").Value = 875
$ws.Range("T23").Value = 16

# Row 24
$ws.Range("D24").Value = 44274
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 130
$ws.Range("N24").Value = 12000
$ws.Range("O24").Value = 12000
$ws.Range("P24").Value = 12000
$ws.Range("Q24").Value = '$/caja 16 kilos'
$ws.Range("S24").Value = 750
$ws.Range("T24").Value = 16

# Row 25
$ws.Range("D25").Value = 44291
$ws.Range("L25").Value = 'Extra (doble especial)'
$ws.Range("M25").Value = 250
$ws.Range("N25").Value = 18000
$ws.Range("O25").Value = 18000
$ws.Range("P25").Value = 18000
$ws.Range("S25").Value = 1000

# Row 26
$ws.Range("D26").Value = 44222
$ws.Range("L26").Value = 'Primera'
$ws.Range("M26").Value = 100
$ws.Range("N26").Value = 18000
$ws.Range("O26").Value = 18000
$ws.Range("P26").Value = 18000
$ws.Range("Q26").Value = '$/caja 16 kilos'
$ws.Range("R26").Value = 'Provincia de Limarí'
$ws.Range("S26").Value = 1125
$ws.Range("T26").Value = 16

# Row 27
$ws.Range("D27").Value = 44301
$ws.Range("M27").Value = 100
$ws.Range("N27").Value = 16000
$ws.Range("O27").Value = 16000
$ws.Range("P27").Value = 16000
$ws.Range("S27").Value = 889

# Row 28
$ws.Range("D28").Value = 44273
$ws.Range("M28").Value = 40
$ws.Range("Q28").Value = '$/caja 16 kilos'
$ws.Range("R28").Value = 'Provincia de Melipilla'
$ws.Range("S28").Value = 938
$ws.Range("T28").Value = 16

# Row 29
$ws.Range("D29").Value = 44273
$ws.Range("L29").Value = 'Primera'
$ws.Range("M29").Value = 50
$ws.Range("N29").Value = 13000
$ws.Range("O29").Value = 13000
$ws.Range("P29").Value = 13000
$ws.Range("Q29").Value = '$/caja 16 kilos'
$ws.Range("S29").Value = 812
$ws.Range("T29").Value = 16

# Row 30
$ws.Range("D30").Value = 44273
$ws.Range("L30").Value = 'Segunda'
$ws.Range("N30").Value = 10000
$ws.Range("O30").Value = 10000
$ws.Range("P30").Value = 10000
$ws.Range("Q30").Value = '$/caja 16 kilos'
$ws.Range("S30").Value = 625
$ws.Range("T30").Value = 16

# Row 31
$ws.Range("D31").Value = 44315
$ws.Range("L31").Value = 'Especial'
$ws.Range("M31").Value = 50
$ws.Range("N31").Value = 24000
$ws.Range("O31").Value = 24000
$ws.Range("P31").Value = 24000
$ws.Range("S31").Value = 1333

# Row 32
$ws.Range("D32").Value = 44315
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 50
$ws.Range("N32").Value = 20000
$ws.Range("O32").Value = 20000
$ws.Range("P32").Value = 20000
$ws.Range("S32").Value = 1111

# Row 33
$ws.Range("D33").Value = 44279
$ws.Range("L33").Value = 'Especial'
$ws.Range("M33").Value = 50
$ws.Range("N33").Value = 14000
$ws.Range("O33").Value = 14000
$ws.Range("P33").Value = 14000
$ws.Range("S33").Value = 778

# Row 34
$ws.Range("D34").Value = 44279
$ws.Range("M34").Value = 100
